# Word TODO list update:
#  - add seven new TODO bullets after "Preserve the state ..."
#  - leave a trailing empty paragraph at the end of the body, like the source
#  - append " when saving to database" before the trailing period of the
#    "Preserve the state ..." bullet (kept as its own run, like the author's
#    original edit)

$d = $word.ActiveDocument

# --- 1. insert the seven new bullet paragraphs -----------------------------
$newBullets = @(
    "Update categories list after restoring from database.",
    "Add vibration and sounds to button presses and other events.",
    "Split categories in the word card to separate entities.",
    "Check for word after OCR to validate if it is found in the dictionary before getting its definition.",
    "Save “first run” variable in preferences and don’t show “how to” every launch.",
    "Callback events like selecting categories or renaming categories can happen twice.",
    "Move key check to splash screen."
)

# Collapse to the end of the "Preserve the state ..." paragraph (which still
# contains the bookmark) and insert one blank paragraph per bullet right
# after it, in order.
$insertionPoint = $d.Paragraphs(3).Range
$insertionPoint.Collapse(0)
for ($i = 0; $i -lt $newBullets.Count; $i++) {
    $insertionPoint.InsertParagraphAfter()
}

for ($i = 0; $i -lt $newBullets.Count; $i++) {
    $d.Paragraphs(4 + $i).Range.Text = $newBullets[$i]
}

# --- 2. trailing empty paragraph at the end of the body --------------------
$lastBullet = $d.Paragraphs(4 + $newBullets.Count - 1).Range
$lastBullet.Collapse(0)
$lastBullet.InsertParagraphAfter()

# --- 3. extend "Preserve the state ..." paragraph --------------------------
# Find the existing sentence (without trailing period) and collapse the
# match range to its end so the new text is inserted right before the ".".
# This must run last: toggling Bold (used below, purely to force the
# inserted text to keep its own run instead of being merged back into its
# neighbour on save) would otherwise bleed into the direct formatting of
# every paragraph inserted afterwards.
$r = $d.Content
$r.Find.Execute("Preserve the state if filtration and sorting on the words list", `
                 $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(" when saving to database")
$r.Bold = 1
$r.Bold = 0
